# fix(publipostage): Try to solve Excel emoji problem
#
# Replace the four "statut" emoji codes used in column A with plain-text /
# alternate-emoji equivalents:
#   📕 -> -3
#   📘 -> ⚠️
#   📙 -> +3
#   📗 -> ✅

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📙" = "+3"
    "📗" = "✅"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        # Force text type so "-3" / "+3" are not auto-coerced into numbers by
        # Excel's value-entry parser, then restore the cell's original
        # (General) style so no visible formatting changes are introduced.
        $cell.NumberFormat = "@"
        $cell.Value2 = $map[$val]
        $cell.Style = "Normal"
    }
}
